$wb = $excel.ActiveWorkbook

# Row 43 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 6185773
$ws.Range("I43").Value = 33933.668
$ws.Range("J43").Value = 9261693
$ws.Range("K43").Value = 33933.668
$ws.Range("L43").Value = 9261693
$ws.Range("M43").Value = -33864.668
$ws.Range("N43").Value = -9261831

# Row 88 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 2473535.8
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 2473535.8
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 2473535.8
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -2474347.8

# Row 91 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 2473535.8
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 2473535.8
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 2473535.8
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -2476343.8

# Row 125 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 2716.5454
$ws.Range("I125").Value = 2047
$ws.Range("K125").Value = 18423
$ws.Range("M125").Value = -15963

# Row 135 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 37038280
$ws.Range("I135").Value = 387.04544
$ws.Range("J135").Value = 200005000
$ws.Range("K135").Value = 3483.40896
$ws.Range("L135").Value = 1800045000
$ws.Range("M135").Value = -948.4089599999998
$ws.Range("N135").Value = -1800050070

# Row 137 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1361.5238
$ws.Range("I137").Value = 929.6799999999999
$ws.Range("K137").Value = 2789.04
$ws.Range("M137").Value = -239.04

# Row 141 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3408.5
$ws.Range("I141").Value = 3637.75
$ws.Range("J141").Value = 2950
$ws.Range("K141").Value = 10913.25
$ws.Range("L141").Value = 8850
$ws.Range("M141").Value = -5733.25
$ws.Range("N141").Value = -19210

# Row 98 on sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H98").Value = 27990
$ws.Range("J98").Value = 27990
$ws.Range("L98").Value = 27990
$ws.Range("N98").Value = -33980

# Row 122 on sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1499.3243
$ws.Range("I122").Value = 1425.0834
$ws.Range("J122").Value = 1636.3846
$ws.Range("K122").Value = 4275.2502
$ws.Range("L122").Value = 4909.1538
$ws.Range("M122").Value = -1825.2502
$ws.Range("N122").Value = -9809.1538

# Row 132 on sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1668.1666
$ws.Range("I132").Value = 1131.1666
$ws.Range("K132").Value = 3393.4998
$ws.Range("M132").Value = -863.4998000000001

# Row 86 on sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4051.375
$ws.Range("I86").Value = 4180.8
$ws.Range("J86").Value = 3835.6667
$ws.Range("K86").Value = 4180.8
$ws.Range("L86").Value = 3835.6667
$ws.Range("M86").Value = -3057.8
$ws.Range("N86").Value = -6081.6667

# Row 89 on sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 4051.375
$ws.Range("I89").Value = 4180.8
$ws.Range("J89").Value = 3835.6667
$ws.Range("K89").Value = 20904
$ws.Range("L89").Value = 19178.3335
$ws.Range("M89").Value = -15288
$ws.Range("N89").Value = -30410.3335

# Row 99 on sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 35715364
$ws.Range("I99").Value = 41667732
$ws.Range("J99").Value = 1152.75
$ws.Range("K99").Value = 41667732
$ws.Range("L99").Value = 1152.75
$ws.Range("M99").Value = -41666234
$ws.Range("N99").Value = -4148.75

# Row 107 on sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1246.3889
$ws.Range("I107").Value = 819
$ws.Range("J107").Value = 1918
$ws.Range("K107").Value = 819
$ws.Range("L107").Value = 1918
$ws.Range("M107").Value = 1101
$ws.Range("N107").Value = -5758

# Row 99 on sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1671.3
$ws.Range("I99").Value = 1641.8334
$ws.Range("K99").Value = 1641.8334
$ws.Range("M99").Value = -143.8334

# Row 120 on sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H120").Value = 18800
$ws.Range("J120").Value = 18800
$ws.Range("L120").Value = 18800
$ws.Range("N120").Value = -26058

# Row 126 on sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1671.3
$ws.Range("I126").Value = 1641.8334
$ws.Range("K126").Value = 4925.5002
$ws.Range("M126").Value = -2455.5002

# Row 132 on sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1744.3334
$ws.Range("I132").Value = 1138.9166
$ws.Range("J132").Value = 4166
$ws.Range("K132").Value = 3416.7498
$ws.Range("L132").Value = 12498
$ws.Range("M132").Value = -886.7498000000001
$ws.Range("N132").Value = -17558

# Row 131 on sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 11112177
$ws.Range("J131").Value = 1133.9147
$ws.Range("L131").Value = 3401.7441
$ws.Range("N131").Value = -13481.7441

# Row 132 on sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 950.6316
$ws.Range("I132").Value = 883.0714
$ws.Range("J132").Value = 1139.8
$ws.Range("K132").Value = 7947.6426
$ws.Range("L132").Value = 10258.2
$ws.Range("M132").Value = -5417.6426
$ws.Range("N132").Value = -15318.2

# Row 139 on sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 1754.5143
$ws.Range("I139").Value = 1791.381
$ws.Range("K139").Value = 5374.143
$ws.Range("M139").Value = -234.143

# Row 102 on sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1101.8889
$ws.Range("I102").Value = 865.1539
$ws.Range("J102").Value = 1717.4
$ws.Range("K102").Value = 865.1539
$ws.Range("L102").Value = 1717.4
$ws.Range("M102").Value = 756.8461
$ws.Range("N102").Value = -4961.4

# Row 107 on sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 627.63635
$ws.Range("I107").Value = 675.5
$ws.Range("J107").Value = 600.2857
$ws.Range("K107").Value = 675.5
$ws.Range("L107").Value = 600.2857
$ws.Range("M107").Value = 1244.5
$ws.Range("N107").Value = -4440.2857

# Row 122 on sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2654.56
$ws.Range("I122").Value = 3093.9285
$ws.Range("J122").Value = 2095.3635
$ws.Range("K122").Value = 9281.7855
$ws.Range("L122").Value = 6286.0905
$ws.Range("M122").Value = -6831.7855
$ws.Range("N122").Value = -11186.0905

# Row 132 on sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1902.3256
$ws.Range("I132").Value = 1794.4667
$ws.Range("J132").Value = 2151.2307
$ws.Range("K132").Value = 5383.4001
$ws.Range("L132").Value = 6453.6921
$ws.Range("M132").Value = -2853.4001
$ws.Range("N132").Value = -11513.6921

# Row 61 on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 821.4
$ws.Range("I61").Value = 790.44446
$ws.Range("K61").Value = 790.44446
$ws.Range("M61").Value = -588.44446

# Row 103 on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H103").Value = 16999.75
$ws.Range("J103").Value = 16999.75
$ws.Range("L103").Value = 16999.75
$ws.Range("N103").Value = -19343.75

# Row 113 on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 821.4
$ws.Range("I113").Value = 790.44446
$ws.Range("K113").Value = 790.44446
$ws.Range("M113").Value = 1379.55554

# Row 122 on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 35716200
$ws.Range("I122").Value = 125001500
$ws.Range("J122").Value = 2079.8
$ws.Range("K122").Value = 375004500
$ws.Range("L122").Value = 6239.400000000001
$ws.Range("M122").Value = -375002050
$ws.Range("N122").Value = -11139.4

# Row 132 on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2390.0571
$ws.Range("I132").Value = 2060.5833
$ws.Range("J132").Value = 3108.9092
$ws.Range("K132").Value = 6181.749899999999
$ws.Range("L132").Value = 9326.7276
$ws.Range("M132").Value = -3651.749899999999
$ws.Range("N132").Value = -14386.7276

# Row 122 on sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 13159370
$ws.Range("I122").Value = 16668269
$ws.Range("J122").Value = 999.75
$ws.Range("K122").Value = 50004807
$ws.Range("L122").Value = 2999.25
$ws.Range("M122").Value = -50002357
$ws.Range("N122").Value = -7899.25

